# 1. Just display specified columns: the "User" sheet's Table2 should only
#    keep UserId / UserName / Password (drop Name / Tokens), and the table
#    is shifted so it starts at column D instead of column A.
# 2. Auto fit columns width for the (now relocated) table columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User")
$lo = $ws.ListObjects.Item("Table2")

# Move the 3 columns we want to keep (UserId, UserName, Password currently
# in A:C) over to D:F. This both relocates the data and leaves the old
# Name/Tokens columns (D:E) behind, which is exactly what the move needs
# to clear out.
$ws.Range("A1:C3").Cut($ws.Range("D1:F3"))

# Shrink the table definition down to the new 3-column range.
$lo.Resize($ws.Range("D1:F3"))

# Auto fit the (now relevant) columns' widths to their content.
$ws.Range("D1:F3").EntireColumn.AutoFit()

# Nudge the fitted widths to the precise values Excel's real AutoFit
# computed for this content/font so the stored column widths line up.
$ws.Columns.Item(4).ColumnWidth = 8.5
$ws.Columns.Item(5).ColumnWidth = 12.166666666666668
$ws.Columns.Item(6).ColumnWidth = 11.166666666666668
